$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "2：用户模式" paragraph: extend the final sentence with three more clauses
#    ("...可以续借，可以看见自己的应还书时间") and keep the trailing "。" but turn
#    it into its own run, with the document's "_GoBack" bookmark sitting
#    between the new text and that final "。" run (matching the target XML,
#    where bookmarkStart/bookmarkEnd moves from the very end of the document
#    to right after "...应还书时间").
#
#    We insert a short, unique marker token immediately before the "。" so we
#    can relocate that exact boundary afterwards (Find.Execute on a Range
#    collapses/repositions that Range onto the match, mirroring real Word).
#    The marker is deleted again once the bookmark is in place, leaving only
#    the two runs + bookmark that the diff calls for.
# ---------------------------------------------------------------------------
$oldSentence = "可以借书，可以还书。"
$marker = "JJJMARKERJJJ"
$newSentence = "可以借书，可以还书，可以续借，可以看见自己的应还书时间" + $marker + "。"

$replaced = $d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                                     $true, 1, $false, $newSentence, 2)

if ($replaced) {
    # Locate the marker; the range collapses onto it.
    $rng = $d.Content
    $rng.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    # Drop a bookmark collapsed right before the marker (i.e. right after
    # "...应还书时间"). Word bookmark names are unique, so (re)adding
    # "_GoBack" here automatically removes it from wherever it used to be
    # (originally at the very end of the document).
    $bmRng = $rng.Duplicate
    $bmRng.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $bmRng)

    # Remove the marker text itself, leaving the bookmark immediately before
    # a standalone "。" run.
    $rng.Text = ""
}

# ---------------------------------------------------------------------------
# 2) Mark the built-in "Normal Table" style as a Quick Style, i.e. add
#    <w:qFormat/> to its style definition in styles.xml.
# ---------------------------------------------------------------------------
$normalTable = $d.Styles.Item("Normal Table")
$normalTable.QuickStyle = $true
